# The deck's Design/theme was switched from the custom "Integral" (Red
# Violet) look to the default Office theme palette. Concretely this means
# the slide master's 12 theme colors (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) change from the Red Violet values to the standard
# Office values - everything else about the theme (fonts, fill/line/
# effect format scheme) is identical between the two themes, so only the
# color scheme needs to change.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.ColorScheme

# RGB() isn't available in this host, so the target sRGB hex values are
# passed through as packed 0xBBGGRR integers (R + G*256 + B*65536), which
# is what PowerPoint's ColorFormat.RGB expects/returns.

$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
